$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q4" sheet right after "总计" (so the tab order
#    becomes 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# NOTE: fetch this *after* the Add() above -- inserting a sheet shifts the
# positional index of every sheet after it, and a Worksheet/Range reference
# obtained beforehand goes stale (its Copy() silently becomes a no-op).
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Bring over the header row (and its styling) from the 2022-Q3 sheet, since
# every quarter sheet shares the exact same column headers/formatting.
# (Column A has no header -- leave A1 untouched/empty, same as every other
# quarter sheet.)
$q3Sheet.Range("B1:H1").Copy($q4.Range("B1:H1"))

# Bring over the styling used for the little numeric index column (A2:A3)
# from the 2022-Q3 sheet as well.
$q3Sheet.Range("A2:A3").Copy($q4.Range("A2:A3"))

# Index column values (0-based row counter).
$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1

# Fund holding data for 2022-Q4. Columns B:G are stored as *text* in the
# source data (even though several look numeric), column H is numeric.
$q4.Range("B2:G3").NumberFormat = "@"

$q4.Range("B2").Value = "160645"
$q4.Range("C2").Value = "鹏华精选回报三年定期开放混合"
$q4.Range("D2").Value = "0.79"
$q4.Range("E2").Value = "74.91"
$q4.Range("F2").Value = "2.82"
$q4.Range("G2").Value = "0.0223"
$q4.Range("H2").Value = 8

$q4.Range("B3").Value = "003670"
$q4.Range("C3").Value = "中融物联网主题灵活配置混合"
$q4.Range("D3").Value = "0.13"
$q4.Range("E3").Value = "92.35"
$q4.Range("F3").Value = "4.29"
$q4.Range("G3").Value = "0.0056"
$q4.Range("H3").Value = 2

# Drop the quote-prefix/text-format styling we applied above so the cells
# come back out with the plain default style (matching the other sheets).
$q4.Range("B2:G3").Style = "Normal"

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q4 row at the top of
#    the data (row 2) and push the existing quarters down by one row.
#    (Values are written literally -- reading `.Value` back out in this
#    host returns the property descriptor, not the cell contents, so we
#    avoid cell-to-cell "copy by reading .Value" entirely.)
# ---------------------------------------------------------------------------

# Copy the index-column (A) styling down into the newly-used row 5 first.
$totalSheet.Range("A4").Copy($totalSheet.Range("A5"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.03

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 7
$totalSheet.Range("D3").Value = 0.96

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q2"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2022-Q1"
$totalSheet.Range("C5").Value = 2
$totalSheet.Range("D5").Value = 0.18

# ---------------------------------------------------------------------------
# 3) Keep the originally-active tab ("2022-Q1") selected -- adding the new
#    sheet above would otherwise steal tab focus for itself.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
